# 지하철노선.xlsx -- "Add files via upload" edit
#
# Business change: the last stop on line 1 ("1호선") is renamed from
# "도봉산" to "잠실새내" and its y-coordinate (column C) becomes 470
# (was 500). Line 2 ("2호선") already lists "잠실새내" as its last
# stop; its B11/C11 cells used to be formulas that mirrored line 1's
# row 11 (`='1호선'!B11`, `='1호선'!C11-30`) -- now that line 1's row
# 11 is itself "잠실새내", those formulas are replaced with the plain
# values they used to evaluate to (500 and 470), so the row no longer
# depends on sheet 1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("1호선")
$ws2 = $wb.Worksheets.Item("2호선")

# --- 1호선 (sheet1): rename last station, update its y-coordinate ---
$ws1.Range("A11").Value = "잠실새내"
$ws1.Range("C11").Value = 470

# --- 2호선 (sheet2): drop the now-circular formulas, keep the values ---
$ws2.Range("B11").Value = 500
$ws2.Range("C11").Value = 470

# --- restore the on-screen selections recorded in the saved file ---
$ws1.Range("A11:C11").Select()
$ws2.Range("E9").Select()
